# "added new pages with new test cases"
#
# - Login sheet gains a new (blank/formatted) second row, and is no longer
#   the active tab.
# - homepagetestdata sheet gains a new test-case value in A2 ("3f2825e0033f")
#   and becomes the active tab/selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Login")
$ws2 = $wb.Worksheets.Item("homepagetestdata")

# Login ("Sheet 1"): new blank row 2 (formatting only, no cell values)
$ws1.Rows.Item(2).RowHeight = 14.05

# homepagetestdata ("Sheet 2"): new test case row with a value
$ws2.Range("A2").Value = "3f2825e0033f"
$ws2.Rows.Item(2).RowHeight = 14.05

# homepagetestdata becomes the active sheet/tab, with A2 selected
$ws2.Activate() | Out-Null
$ws2.Range("A2").Select() | Out-Null
